$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing Row 2 ---
$ws.Cells.Item(2, 3).Value = 'OrderedDict([(''id'', 40), (''password'', ''pbkdf2_sha256$390000$MUIB8FwHiRDV6Phpz2UsCA$utHSfT/1pqI47SdblTUanKF+KhLFtCOwuLd7kwbrKWA=''), (''last_login'', ''2022-12-07T14:17:40.845462+05:45''), (''is_superuser'', False), (''first_name'', ''private user''), (''last_name'', ''''), (''is_staff'', False), (''is_active'', True), (''date_joined'', ''2022-11-29T11:35:16.469786+05:45''), (''sex'', None), (''phone'', ''9808098''), (''permanent_address'', None), (''email'', ''manojdas.py@gmail.com''), (''username'', ''user''), (''state_name'', ''Madhesh Province''), (''district_name'', ''Dhanusha''), (''is_verified'', True), (''is_applyForVerified'', True), (''apply_role_type'', 3), (''company_name'', ''private company''), (''union_type'', ''district''), (''union_name'', ''central@gmail.com''), (''approved_name'', None), (''approved_email'', ''manojdas.py@gmail.com''), (''approved_signature'', None), (''approved_company_name'', ''md company''), (''image'', ''/media/user/profile/my_photo.JPG''), (''signature'', ''/media/user/documents/signature_ceo.jpeg''), (''member_admin'', ''admin''), (''created_at'', ''2022-11-29T11:35:16.940528+05:45''), (''updated_at'', ''2022-12-04T12:07:18.398027+05:45''), (''role'', 5), (''groups'', []), (''user_permissions'', [])])'
$ws.Cells.Item(2, 4).Value = 'OrderedDict([(''id'', 23), (''owner_full_name'', ''private user''), (''company_name'', ''private company''), (''union_type'', ''district''), (''union_name'', ''central@gmail.com''), (''business_name'', ''plywood_industry''), (''business_price_category'', ''old_unregistered_non_payment''), (''is_reniew'', False), (''payment_rupees'', 20000), (''voucher_number'', ''None''), (''municipality'', ''Mithila''), (''ward_number'', ''2''), (''tole'', ''topa''), (''authorized_capital_of_company'', 98), (''paid_up_capital_of_company'', 98), (''transaction_amount'', 456), (''tax_amount'', 987), (''number_of_employees'', 98), (''approved_name'', ''central''), (''approved_email'', ''central@gmail.com''), (''approved_signature'', ''/media/user/documents/signature_eLOwzGt.jpeg''), (''approved_company_name'', None), (''certificate_citizenship'', ''/media/user/documents/Doc1_srdTaWL.jpg''), (''certificate_company_registration'', ''/media/user/documents/Doc2_mG3bf1z.jpg''), (''provisional_account_number'', ''/media/user/documents/Doc3_VPLgpYj.jpg''), (''auditing'', None), (''tax_paid_certificate'', ''/media/user/documents/Doc1_NohrdG2.jpg''), (''industry_registration'', None), (''income_certificate'', None), (''member_approved_certificate'', None), (''self_declaration_not_on_the_black_list'', None), (''bank_voucher'', ''/media/user/documents/qr-new.jpg''), (''industry_certificate_front_side'', ''/media/user/documents/my_photo.JPG''), (''industry_certificate_back_side'', ''/media/user/documents/my_photo_hGyyNkP.JPG''), (''application_certificate'', ''/media/user/documents/PDF_SIGNING_USING_DIGITAL_SIGNATURE_-_Google_Docs_4.pdf''), (''physical_year'', ''2078''), (''created_at'', ''2022-11-29T11:40:29.785957+05:45''), (''updated_at'', ''2022-12-09T13:53:51.136588+05:45''), (''user'', 40)])'
$ws.Cells.Item(2, 5).Value = 'central_ceo'
$ws.Cells.Item(2, 12).Value = '2022-12-09T13:53:44.554927+05:45'
$ws.Cells.Item(2, 19).Value = 38

# --- Add new Row 3 ---
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(3, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = 24
$ws.Cells.Item(3, 3).Value = 'OrderedDict([(''id'', 43), (''password'', ''pbkdf2_sha256$390000$N7ybMyN8rPjN2PAqbTyebZ$6bP1q2HSE9ne7oSQXLmesj8afNFthAjpKLKBz94DVm8=''), (''last_login'', ''2022-12-04T16:09:10.190495+05:45''), (''is_superuser'', False), (''first_name'', ''मनोज दास''), (''last_name'', ''''), (''is_staff'', False), (''is_active'', True), (''date_joined'', ''2022-12-04T16:09:09.463723+05:45''), (''sex'', ''male''), (''phone'', ''98080''), (''permanent_address'', None), (''email'', ''a@gmail.com''), (''username'', ''a''), (''state_name'', ''Province No. 1''), (''district_name'', ''Ilam''), (''is_verified'', True), (''is_applyForVerified'', True), (''apply_role_type'', 5), (''company_name'', ''asd''), (''union_type'', ''district''), (''union_name'', ''central@gmail.com''), (''approved_name'', None), (''approved_email'', None), (''approved_signature'', None), (''approved_company_name'', None), (''image'', None), (''signature'', None), (''member_admin'', ''member''), (''created_at'', ''2022-12-04T16:09:09.750044+05:45''), (''updated_at'', ''2022-12-04T17:16:38.929999+05:45''), (''role'', 5), (''groups'', []), (''user_permissions'', [])])'
$ws.Cells.Item(3, 4).Value = 'OrderedDict([(''id'', 24), (''owner_full_name'', ''a''), (''company_name'', ''asd''), (''union_type'', ''district''), (''union_name'', ''central@gmail.com''), (''business_name'', ''plywood_industry''), (''business_price_category'', ''less_than_2_crode''), (''is_reniew'', False), (''payment_rupees'', 10000), (''voucher_number'', ''None''), (''municipality'', ''Mangsebung''), (''ward_number'', ''8''), (''tole'', ''3''), (''authorized_capital_of_company'', 87), (''paid_up_capital_of_company'', 87), (''transaction_amount'', 8797), (''tax_amount'', 98), (''number_of_employees'', 890), (''approved_name'', None), (''approved_email'', None), (''approved_signature'', None), (''approved_company_name'', None), (''certificate_citizenship'', None), (''certificate_company_registration'', None), (''provisional_account_number'', None), (''auditing'', None), (''tax_paid_certificate'', None), (''industry_registration'', None), (''income_certificate'', None), (''member_approved_certificate'', None), (''self_declaration_not_on_the_black_list'', None), (''bank_voucher'', None), (''industry_certificate_front_side'', None), (''industry_certificate_back_side'', None), (''application_certificate'', ''/media/user/documents/71-72.pdf''), (''physical_year'', ''2079''), (''created_at'', ''2022-12-04T16:09:53.993304+05:45''), (''updated_at'', ''2022-12-07T17:44:22.720251+05:45''), (''user'', 43)])'
$ws.Cells.Item(3, 5).Value = 'central_ceo'
$ws.Cells.Item(3, 6).Value = $false
$ws.Cells.Item(3, 7).Value = $false
$ws.Cells.Item(3, 8).Value = $false
$ws.Cells.Item(3, 9).Value = $false
$ws.Cells.Item(3, 11).Value = '2022-12-04T16:09:54.158928+05:45'
$ws.Cells.Item(3, 12).Value = '2022-12-06T17:31:40.564555+05:45'
$ws.Cells.Item(3, 13).Value = $true

# Force text (not numeric) storage for N3/O3/P3, then clear the format so no extra style sticks
$cellN3 = $ws.Cells.Item(3, 14)
$cellN3.NumberFormat = "@"
$cellN3.Value = '234'
$cellN3.Style = "Normal"

$cellO3 = $ws.Cells.Item(3, 15)
$cellO3.NumberFormat = "@"
$cellO3.Value = '3234'
$cellO3.Style = "Normal"

$cellP3 = $ws.Cells.Item(3, 16)
$cellP3.NumberFormat = "@"
$cellP3.Value = '1'
$cellP3.Style = "Normal"

$ws.Cells.Item(3, 19).Value = 38
